$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new publication row (row 24) documenting a manuscript under review.
$ws.Range("A24").Value = "prep"
$ws.Range("C24").Value = "(Under Review)"
$ws.Range("D24").Value = "Lead Exposure and Antisocial Behavior: A Systematic Review Protocol"
$ws.Range("E24").Value = "Under review at Environment International"
$ws.Range("B24").Value = "EF Kirrane, RM Shaffer, J Forsyth, G Ferraro, C Hill, L Carlson, K Hester, CC Lanfear, H Hu"
